$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2045.9231
$ws.Range("I11").Value = 2045.9231
$ws.Range("K11").Value = 2045.9231
$ws.Range("M11").Value = -1905.9231

$ws.Range("H18").Value = 166683420
$ws.Range("I18").Value = 250005120
$ws.Range("K18").Value = 250005120
$ws.Range("M18").Value = -250004836

$ws.Range("H55").Value = 185.19048
$ws.Range("I55").Value = 189.95
$ws.Range("K55").Value = 189.95
$ws.Range("M55").Value = 24.05000000000001

$ws.Range("H98").Value = 38464596
$ws.Range("I98").Value = 40002980
$ws.Range("K98").Value = 40002980
$ws.Range("M98").Value = -40001482

$ws.Range("H113").Value = 73244260
$ws.Range("J113").Value = 93765576
$ws.Range("L113").Value = 93765576
$ws.Range("N113").Value = -93772084

$ws.Range("H122").Value = 38464596
$ws.Range("I122").Value = 40002980
$ws.Range("K122").Value = 120008940
$ws.Range("M122").Value = -120006490

$ws.Range("H129").Value = 1220.7778
$ws.Range("I129").Value = 585.7
$ws.Range("K129").Value = 1757.1
$ws.Range("M129").Value = 3242.9

$ws.Range("H137").Value = 4177.1055
$ws.Range("I137").Value = 2799.6086
$ws.Range("K137").Value = 8398.825800000001
$ws.Range("M137").Value = -5848.825800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 11375
$ws.Range("I26").Value = 11375
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 11375
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -11045

$ws.Range("H32").Value = 1926296.2
$ws.Range("I32").Value = 1926296.2
$ws.Range("K32").Value = 1926296.2
$ws.Range("M32").Value = -1926009.2

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H45").Value = 2264.4443
$ws.Range("I45").Value = 1334.25
$ws.Range("K45").Value = 1334.25
$ws.Range("M45").Value = -957.25

$ws.Range("H61").Value = 33338944
$ws.Range("I61").Value = 1578.1578
$ws.Range("K61").Value = 1578.1578
$ws.Range("M61").Value = -1366.1578

$ws.Range("H132").Value = 4982.0156
$ws.Range("I132").Value = 2884.425
$ws.Range("J132").Value = 8478
$ws.Range("K132").Value = 8653.275000000001
$ws.Range("L132").Value = 25434
$ws.Range("M132").Value = -6123.275000000001
$ws.Range("N132").Value = -30494

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws.Range("H136").Value = 33338944
$ws.Range("I136").Value = 1578.1578
$ws.Range("K136").Value = 4734.4734
$ws.Range("M136").Value = -2184.4734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1386.6875
$ws.Range("I94").Value = 932.9167
$ws.Range("J94").Value = 2748
$ws.Range("K94").Value = 932.9167
$ws.Range("L94").Value = 2748
$ws.Range("M94").Value = -481.9167
$ws.Range("N94").Value = -3650

$ws.Range("H134").Value = 5323682
$ws.Range("I134").Value = 8066321.5
$ws.Range("K134").Value = 24198964.5
$ws.Range("M134").Value = -24196429.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5412.2812
$ws.Range("I31").Value = 1957.3667
$ws.Range("J31").Value = 8460.735000000001
$ws.Range("K31").Value = 1957.3667
$ws.Range("L31").Value = 8460.735000000001
$ws.Range("M31").Value = -1662.3667
$ws.Range("N31").Value = -9050.735000000001

$ws.Range("H32").Value = 2755.25
$ws.Range("I32").Value = 1755
$ws.Range("K32").Value = 1755
$ws.Range("M32").Value = -1439

$ws.Range("H34").Value = 5412.2812
$ws.Range("I34").Value = 1957.3667
$ws.Range("J34").Value = 8460.735000000001
$ws.Range("K34").Value = 1957.3667
$ws.Range("L34").Value = 8460.735000000001
$ws.Range("M34").Value = -1755.3667
$ws.Range("N34").Value = -8864.735000000001

$ws.Range("H122").Value = 44993.918
$ws.Range("I122").Value = 2569.1428
$ws.Range("K122").Value = 7707.428400000001
$ws.Range("M122").Value = -5257.428400000001

$ws.Range("H132").Value = 8265.412
$ws.Range("I132").Value = 4303.4
$ws.Range("K132").Value = 12910.2
$ws.Range("M132").Value = -10380.2

$ws.Range("H133").Value = 33053.715
$ws.Range("J133").Value = 33053.715
$ws.Range("L133").Value = 33053.715
$ws.Range("N133").Value = -38113.715

$ws.Range("H134").Value = 6799.857
$ws.Range("I134").Value = 3017
$ws.Range("J134").Value = 8313
$ws.Range("K134").Value = 9051
$ws.Range("L134").Value = 24939
$ws.Range("M134").Value = -6516
$ws.Range("N134").Value = -30009

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2778395.2
$ws.Range("I12").Value = 2163.4
$ws.Range("J12").Value = 3846176.8
$ws.Range("K12").Value = 6490.200000000001
$ws.Range("L12").Value = 11538530.4
$ws.Range("M12").Value = -6317.200000000001
$ws.Range("N12").Value = -11538876.4

$ws.Range("H44").Value = 570.8333
$ws.Range("I44").Value = 298.33334
$ws.Range("K44").Value = 895.0000200000001
$ws.Range("M44").Value = -497.0000200000001

$ws.Range("H87").Value = 5599.3335
$ws.Range("I87").Value = 4899.5
$ws.Range("J87").Value = 6999
$ws.Range("K87").Value = 14698.5
$ws.Range("L87").Value = 20997
$ws.Range("M87").Value = -13450.5
$ws.Range("N87").Value = -23493

$ws.Range("H90").Value = 5599.3335
$ws.Range("I90").Value = 4899.5
$ws.Range("J90").Value = 6999
$ws.Range("K90").Value = 44095.5
$ws.Range("L90").Value = 62991
$ws.Range("M90").Value = -37855.5
$ws.Range("N90").Value = -75471

$ws.Range("H98").Value = 55557710
$ws.Range("J98").Value = 100003760
$ws.Range("L98").Value = 300011280
$ws.Range("N98").Value = -300014276

$ws.Range("H132").Value = 6106.8276
$ws.Range("I132").Value = 2134.7058
$ws.Range("J132").Value = 11734
$ws.Range("K132").Value = 19212.3522
$ws.Range("L132").Value = 105606
$ws.Range("M132").Value = -16682.3522
$ws.Range("N132").Value = -110666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("N29").Value = 0

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H80").Value = 5733.3335
$ws.Range("I80").Value = 4852.5
$ws.Range("J80").Value = 7495
$ws.Range("K80").Value = 4852.5
$ws.Range("L80").Value = 7495
$ws.Range("M80").Value = -3854.5
$ws.Range("N80").Value = -9491

$ws.Range("H83").Value = 5733.3335
$ws.Range("I83").Value = 4852.5
$ws.Range("J83").Value = 7495
$ws.Range("K83").Value = 24262.5
$ws.Range("L83").Value = 37475
$ws.Range("M83").Value = -19270.5
$ws.Range("N83").Value = -47459

$ws.Range("H122").Value = 4543514.5
$ws.Range("J122").Value = 5204.923
$ws.Range("L122").Value = 15614.769
$ws.Range("N122").Value = -20514.769

$ws.Range("H132").Value = 3254.5293
$ws.Range("I132").Value = 1923.275
$ws.Range("K132").Value = 5769.825000000001
$ws.Range("M132").Value = -3239.825000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H32").Value = 12506.5
$ws.Range("I32").Value = 12506.5
$ws.Range("K32").Value = 12506.5
$ws.Range("M32").Value = -12189.5

$ws.Range("H46").Value = 3971333.8
$ws.Range("I46").Value = 2699.9333
$ws.Range("J46").Value = 8550527
$ws.Range("K46").Value = 2699.9333
$ws.Range("L46").Value = 8550527
$ws.Range("M46").Value = -2511.9333
$ws.Range("N46").Value = -8550903

$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50450

$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51560

$ws.Range("H88").Value = 61999.5
$ws.Range("J88").Value = 61999.5
$ws.Range("L88").Value = 61999.5
$ws.Range("N88").Value = -62855.5

$ws.Range("H91").Value = 61999.5
$ws.Range("J91").Value = 61999.5
$ws.Range("L91").Value = 61999.5
$ws.Range("N91").Value = -64963.5

$ws.Range("H100").Value = 4208.826
$ws.Range("I100").Value = 1758.3636
$ws.Range("K100").Value = 1758.3636
$ws.Range("M100").Value = -1217.3636

$ws.Range("H132").Value = 9443280
$ws.Range("I132").Value = 20836212
$ws.Range("J132").Value = 14646.034
$ws.Range("K132").Value = 62508636
$ws.Range("L132").Value = 43938.102
$ws.Range("M132").Value = -62506106
$ws.Range("N132").Value = -48998.102

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I5").Value = 9000000
$ws.Range("K5").Value = 9000000
$ws.Range("M5").Value = -8999888

$ws.Range("H51").Value = 7999.6665
$ws.Range("I51").Value = 7999.6665
$ws.Range("K51").Value = 7999.6665
$ws.Range("M51").Value = -7489.6665

$ws.Range("H126").Value = 125003420
$ws.Range("J126").Value = 2902.1667
$ws.Range("L126").Value = 8706.500100000001
$ws.Range("N126").Value = -13646.5001

$ws.Range("H133").Value = 166000
$ws.Range("J133").Value = 166000
$ws.Range("L133").Value = 166000
$ws.Range("N133").Value = -176120
